$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.008.73'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").Value = '1.645.70'
$ws.Range("E3").Value = '  -1.32%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.18'
$ws.Range("E5").Value = '  +2.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5216'
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2606'
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06362'
$ws.Range("E9").Value = '  +0.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.79'
$ws.Range("E10").Value = '  -1.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07682'
$ws.Range("E11").Value = '  +1.97%  '
$ws.Range("D12").Value = '1.651.00'
$ws.Range("E12").Value = '  -1.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.421'
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("E14").Value = '  -1.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5539'
$ws.Range("E15").Value = '  +1.96%  '
$ws.Range("D16").Value = '0.0₅8293'
$ws.Range("E16").Value = '  +3.55%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.73'
$ws.Range("E17").Value = '  -2.47%  '
$ws.Range("D18").Value = '26.016.03'
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.711'
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '188.14'
$ws.Range("E21").Value = '  +0.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.16'
$ws.Range("E22").Value = '  -1.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.252'
$ws.Range("E23").Value = '  +0.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.31'
$ws.Range("E25").Value = '  -3.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1223'
$ws.Range("E26").Value = '  -1.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.403'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.88'
$ws.Range("E28").Value = '  +0.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.388'
$ws.Range("E29").Value = '  +1.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05943'
$ws.Range("E30").Value = '  -5.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.264'
$ws.Range("E31").Value = '  -1.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.396'
$ws.Range("E32").Value = '  -0.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.403'
$ws.Range("E33").Value = '  -2.78%  '
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9948'
$ws.Range("E35").Value = '  -0.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.395'
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.754'
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5633'
$ws.Range("E38").Value = '  -5.96%  '
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.867'
$ws.Range("E40").Value = '  -2.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8530'
$ws.Range("E41").Value = '  -1.08%  '
$ws.Range("E42").Value = '  -0.21%  '
$ws.Range("D43").Value = '1.028.02'
$ws.Range("E43").Value = '  -7.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.81'
$ws.Range("E44").Value = '  -1.82%  '
$ws.Range("D45").Value = '1.795.65'
$ws.Range("E45").Value = '  -1.30%  '
$ws.Range("D46").Value = '0.0₈107'
$ws.Range("E46").Value = '  -0.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.68'
$ws.Range("E47").Value = '  +0.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.002'
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.048'
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05144'
$ws.Range("E50").Value = '  -1.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4217'
$ws.Range("E51").Value = '  -0.44%  '
